# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume update described by the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.872.57'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.891.87'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("D4").Value = '''1.009'
$ws.Range("E4").Value = '  +0.97%  '
$ws.Range("D5").Value = '''324.05'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '''1.008'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("D8").Value = '''0.3813'
$ws.Range("E8").Value = '  -2.29%  '
$ws.Range("D9").Value = '''0.07722'
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("D10").Value = '''0.9648'
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("D11").Value = '''21.93'
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = '1.885.78'
$ws.Range("E12").Value = '  -2.95%  '
$ws.Range("D13").Value = '''6.975'
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("D14").Value = '''5.683'
$ws.Range("E14").Value = '  -2.52%  '
$ws.Range("D15").Value = '''0.07082'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '''1.010'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = '''83.40'
$ws.Range("E17").Value = '  -4.94%  '
$ws.Range("D18").Value = '''0.000009482'
$ws.Range("E18").Value = '  -4.80%  '
$ws.Range("D19").Value = '''16.81'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").Value = '''1.007'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Value = '28.882.68'
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").Value = '''5.381'
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("D23").Value = '''10.93'
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("D24").Value = '2.203.63'
$ws.Range("E24").Value = '  +1.28%  '
$ws.Range("D25").Value = '''2.080'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '''157.31'
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").Value = '''19.10'
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("D28").Value = '''5.647'
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("D29").Value = '''117.29'
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").Value = '''1.825'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").Value = '''0.09298'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").Value = '''0.8578'
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("D33").Value = '''5.083'
$ws.Range("E33").Value = '  -2.59%  '
$ws.Range("D34").Value = '''1.247'
$ws.Range("E34").Value = '  -6.23%  '
$ws.Range("D35").Value = '''3.078'
$ws.Range("E35").Value = '  -1.92%  '
$ws.Range("D36").Value = '''1.154'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").Value = '''0.05649'
$ws.Range("E37").Value = '  -2.38%  '
$ws.Range("D38").Value = '''0.02042'
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").Value = '''0.5529'
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("D40").Value = '''7.440'
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").Value = '''0.1753'
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '''9.236'
$ws.Range("E42").Value = '  -5.34%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '''0.000002845'
$ws.Range("E43").Value = '  -15.76%  '
$ws.Range("D44").Value = '''2.692'
$ws.Range("E44").Value = '  +3.67%  '
$ws.Range("D45").Value = '''0.5184'
$ws.Range("E45").Value = '  -2.98%  '
$ws.Range("D46").Value = '''11.14'
$ws.Range("E46").Value = '  -7.32%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''2.077'
$ws.Range("E47").Value = '  -6.01%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.06768'
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("D49").Value = '''1.775'
$ws.Range("E49").Value = '  -4.03%  '
$ws.Range("D50").Value = '''110.20'
$ws.Range("E50").Value = '  -2.62%  '
$ws.Range("D51").Value = '''0.2956'
$ws.Range("E51").Value = '  -1.46%  '
